$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Cuenta 1"
$ws.Range("C2").Value = 564856

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Cuenta 2"
$ws.Range("C3").Value = 45200
